# Add data for 2022-08-07: update "through" date from 07-29 to 07-30
# and refresh the July row / Total row figures accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the workbook's window title reference
$ws.Name = "Through 2022-07-30"

# Update the "July (through 07-29)" label cell to "July (through 07-30)"
$ws.Range("A8").Value = "July (through 07-30)"

# Update July row (row 8) values
$ws.Range("B8").Value = 37
$ws.Range("C8").Value = 53
$ws.Range("D8").Value = 71
$ws.Range("E8").Value = 69
$ws.Range("F8").Value = 51
$ws.Range("G8").Value = 143
$ws.Range("H8").Value = 142
$ws.Range("I8").Value = 164

# Update Total row (row 9) values
$ws.Range("B9").Value = 162
$ws.Range("C9").Value = 301
$ws.Range("D9").Value = 461
$ws.Range("E9").Value = 422
$ws.Range("F9").Value = 302
$ws.Range("G9").Value = 615
$ws.Range("H9").Value = 902
$ws.Range("I9").Value = 970
